$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Geometry Garage")

$co2 = $ws.ChartObjects().Add(1044.2958984375, 190.12496062992125, 270.4375, 216)
$co2.Name = "Chart 2"
$chart2 = $co2.Chart
$chart2.ChartType = 51

$ser = $chart2.SeriesCollection().NewSeries()
$ser.Formula = "=SERIES('Geometry Garage'!`$L`$4,'Geometry Garage'!`$M`$3:`$P`$3,'Geometry Garage'!`$M`$4:`$P`$4,1)"

Write-Host "Series count:" $chart2.SeriesCollection().Count
$s = $chart2.SeriesCollection(1)
Write-Host "Name:" $s.Name
Write-Host "Formula:" $s.Formula
